# "Ten imputations instead of five" - update the reported AIC-ish metric
# values (column B) for each functional-form row on the active sheet.
# The values are stored as text (shared strings) in the workbook, so we
# force text entry with a leading apostrophe (quote-prefix) and then
# restore the "Normal" style so no stray number-format is left on the
# cells themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
  2  = "182.7"
  3  = "183.9"
  4  = "182.4"
  5  = "184.9"
  6  = "184.1"
  7  = "185.9"
  8  = "187.6"
  9  = "183.8"
  10 = "185.8"
  11 = "185.0"
  12 = "184.1"
  13 = "183.7"
  14 = "185.1"
  15 = "185.3"
}

foreach ($row in $newValues.Keys) {
  $cell = $ws.Range("B$row")
  $cell.Value = "'" + $newValues[$row]
  $cell.Style = "Normal"
}
